$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Make Parts the active/selected sheet (matches activeTab + tabSelected in target)
$ws.Activate()

# Insert a new row at row 14, shifting the existing rows 14-38 down to 15-39
$ws.Rows.Item(14).Insert()

# The inserted row doesn't inherit the bordered "line item" formatting used by
# the rest of the table, so copy it over from the row above (row 13), which
# uses the identical border/style pattern.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
foreach ($col in $cols) {
  $ws.Range($col + "13").Copy()
  $ws.Range($col + "14").PasteSpecial(-4122)
}
[void]($excel.CutCopyMode = 0)

# H14/I14 (the part number / MPN) were entered without the bordered style
$ws.Range("H14:I14").ClearFormats()

# New BOM line: 12V PTC fuse added ahead of item 9 (U1)
$ws.Range("E14").Value = "8a"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "R8"
$ws.Range("H14").Value = "507-1818-1-ND  "
$ws.Range("I14").Value = "0ZCM0010FF2G"
$ws.Range("J14").Value = "FUSE PTC 100MA "

# The printed area grows by one row to keep the whole table on the printout
$ws.PageSetup.PrintArea = "B2:L39"

# Leave the selection where the author left it after adding the row
[void]($ws.Range("J15").Select())
